$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 57428638.29
$ws.Range("P2").Value = 94.5559004638
$ws.Range("Q2").Value = 902958387.25
$ws.Range("R2").Value = 1486.7154425051
$ws.Range("S2").Value = 519190776.5
$ws.Range("T2").Value = 854.8444268618
$ws.Range("U2").Value = -15490.68
$ws.Range("V2").Value = -0.0255053095
$ws.Range("W2").Value = 126287.67
$ws.Range("X2").Value = 0.2079318735
$ws.Range("Y2").Value = 3080684.04
$ws.Range("Z2").Value = 5.0723273635
$ws.Range("AA2").Value = 3313796.2
$ws.Range("AB2").Value = 5.4561451042
$ws.Range("AC2").Value = 60735118.6
$ws.Range("AD2").Value = 76.45199261960001
